# Lernfeld 4 Teil 5 / SYS.2.2.3.A1 Maßnahmen.xlsx
# - Fix wording in B5 (add comma)
# - Add a new "Begründung" row (row 6) explaining the measure
# - Update the active selection to the next empty cell (B7)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing memo text (added comma after "Mitarbeiter")
$ws.Range("B5").Value = "Memo an alle Mitarbeiter, welche Cloud benutzt werden soll"

# Append new row with justification for the measure
$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Zur Vereinheitlichung des Arbeitsprozesses"

# Move the selection to where the next entry would go, as Excel does
# after typing into the last used row
[void]$ws.Range("B7").Select()
